# PV-94: Support for plans without sticky-ids or levels
# Rename header columns to the new terminology:
#   A1: "Unique Sticky ID" -> "Row ID"
#   C1: "Name"             -> "Task"
#   E1: "Start"            -> "Start Date"
#   F1: "Finish"           -> "End Date"
# and move the active selection to F1 (was C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

$ws.Range("F1").Select()
